$p = $ppt.ActivePresentation
$s = $p.Slides.Item(5)
$sh = $s.Shapes.Item(2)
$tf = $sh.TextFrame
$tr = $tf.TextRange

# --- Change 1: paragraph 5 wording tweak ("as same complex" -> "as same complex level") ---
$para5 = $tr.Paragraphs(5)
$run5 = $para5.Runs(1)
$run5.Text = "The test architecture must align with application development architecture. So that test script can be updated as same complex level as application development."

# --- Change 2: paragraph 6, merge the trailing three runs into one ---
$para6 = $tr.Paragraphs(6)
$runA = $para6.Runs(6)
$runB = $para6.Runs(7)
$runC = $para6.Runs(8)
$combinedLen = $runA.Length + $runB.Length + $runC.Length
$mergedRange = $tr.Characters($runA.Start, $combinedLen)
$mergedRange.Text = " data structure in web application development, so the test data structure also should align to Json."
